# Update Top33_DataComp worksheet with refreshed M2/FX length & date statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellUpdates = @(
    @{ Row=2; Col=3; Value=342 }  # C2: 341 -> 342
    @{ Row=2; Col=6; Value=45444 }  # F2: 45413 -> 45444
    @{ Row=2; Col=7; Value=30319 }  # G2: 30256 -> 30319
    @{ Row=2; Col=8; Value=45505 }  # H2: 45446 -> 45505
    @{ Row=3; Col=5; Value=30256 }  # E3: 30225 -> 30256
    @{ Row=3; Col=6; Value=45444 }  # F3: 45413 -> 45444
    @{ Row=4; Col=5; Value=30256 }  # E4: 30225 -> 30256
    @{ Row=4; Col=6; Value=45444 }  # F4: 45413 -> 45444
    @{ Row=4; Col=7; Value=30319 }  # G4: 30256 -> 30319
    @{ Row=4; Col=8; Value=45505 }  # H4: 45446 -> 45505
    @{ Row=5; Col=5; Value=30256 }  # E5: 30225 -> 30256
    @{ Row=5; Col=6; Value=45444 }  # F5: 45413 -> 45444
    @{ Row=5; Col=7; Value=30319 }  # G5: 30256 -> 30319
    @{ Row=5; Col=8; Value=45505 }  # H5: 45446 -> 45505
    @{ Row=6; Col=3; Value=451 }  # C6: 449 -> 451
    @{ Row=6; Col=6; Value=45444 }  # F6: 45383 -> 45444
    @{ Row=6; Col=7; Value=30319 }  # G6: 30256 -> 30319
    @{ Row=6; Col=8; Value=45505 }  # H6: 45446 -> 45505
    @{ Row=7; Col=5; Value=30256 }  # E7: 30195 -> 30256
    @{ Row=7; Col=6; Value=45444 }  # F7: 45383 -> 45444
    @{ Row=7; Col=7; Value=30319 }  # G7: 30256 -> 30319
    @{ Row=7; Col=8; Value=45505 }  # H7: 45446 -> 45505
    @{ Row=8; Col=4; Value=418 }  # D8: 416 -> 418
    @{ Row=8; Col=5; Value=30256 }  # E8: 30195 -> 30256
    @{ Row=8; Col=6; Value=45444 }  # F8: 45383 -> 45444
    @{ Row=8; Col=8; Value=45505 }  # H8: 45446 -> 45505
    @{ Row=9; Col=5; Value=30256 }  # E9: 30225 -> 30256
    @{ Row=9; Col=6; Value=45444 }  # F9: 45413 -> 45444
    @{ Row=9; Col=7; Value=30319 }  # G9: 30256 -> 30319
    @{ Row=9; Col=8; Value=45505 }  # H9: 45446 -> 45505
    @{ Row=10; Col=4; Value=491 }  # D10: 489 -> 491
    @{ Row=10; Col=5; Value=30256 }  # E10: 30225 -> 30256
    @{ Row=10; Col=6; Value=45444 }  # F10: 45413 -> 45444
    @{ Row=10; Col=8; Value=45505 }  # H10: 45446 -> 45505
    @{ Row=11; Col=5; Value=30225 }  # E11: 30195 -> 30225
    @{ Row=11; Col=6; Value=45413 }  # F11: 45383 -> 45413
    @{ Row=11; Col=7; Value=30319 }  # G11: 30256 -> 30319
    @{ Row=11; Col=8; Value=45505 }  # H11: 45446 -> 45505
    @{ Row=12; Col=3; Value=379 }  # C12: 378 -> 379
    @{ Row=12; Col=4; Value=360 }  # D12: 358 -> 360
    @{ Row=12; Col=6; Value=45444 }  # F12: 45413 -> 45444
    @{ Row=12; Col=8; Value=45505 }  # H12: 45446 -> 45505
    @{ Row=13; Col=3; Value=474 }  # C13: 472 -> 474
    @{ Row=13; Col=6; Value=45413 }  # F13: 45352 -> 45413
    @{ Row=13; Col=7; Value=30319 }  # G13: 30256 -> 30319
    @{ Row=13; Col=8; Value=45505 }  # H13: 45446 -> 45505
    @{ Row=14; Col=4; Value=404 }  # D14: 402 -> 404
    @{ Row=14; Col=8; Value=45505 }  # H14: 45446 -> 45505
    @{ Row=15; Col=7; Value=30286 }  # G15: 30225 -> 30286
    @{ Row=15; Col=8; Value=45505 }  # H15: 45446 -> 45505
    @{ Row=16; Col=3; Value=462 }  # C16: 461 -> 462
    @{ Row=16; Col=4; Value=418 }  # D16: 416 -> 418
    @{ Row=16; Col=6; Value=45413 }  # F16: 45383 -> 45413
    @{ Row=16; Col=8; Value=45505 }  # H16: 45446 -> 45505
    @{ Row=17; Col=3; Value=378 }  # C17: 376 -> 378
    @{ Row=17; Col=4; Value=402 }  # D17: 400 -> 402
    @{ Row=17; Col=6; Value=45444 }  # F17: 45383 -> 45444
    @{ Row=17; Col=8; Value=45505 }  # H17: 45446 -> 45505
    @{ Row=18; Col=4; Value=268 }  # D18: 266 -> 268
    @{ Row=18; Col=5; Value=30256 }  # E18: 30195 -> 30256
    @{ Row=18; Col=6; Value=45444 }  # F18: 45383 -> 45444
    @{ Row=18; Col=8; Value=45505 }  # H18: 45446 -> 45505
    @{ Row=19; Col=4; Value=406 }  # D19: 404 -> 406
    @{ Row=19; Col=5; Value=30256 }  # E19: 30225 -> 30256
    @{ Row=19; Col=6; Value=45444 }  # F19: 45413 -> 45444
    @{ Row=19; Col=8; Value=45505 }  # H19: 45446 -> 45505
    @{ Row=20; Col=3; Value=489 }  # C20: 487 -> 489
    @{ Row=20; Col=6; Value=45444 }  # F20: 45383 -> 45444
    @{ Row=20; Col=7; Value=30319 }  # G20: 30256 -> 30319
    @{ Row=20; Col=8; Value=45505 }  # H20: 45446 -> 45505
    @{ Row=21; Col=3; Value=317 }  # C21: 316 -> 317
    @{ Row=21; Col=6; Value=45413 }  # F21: 45383 -> 45413
    @{ Row=21; Col=7; Value=30319 }  # G21: 30256 -> 30319
    @{ Row=21; Col=8; Value=45505 }  # H21: 45446 -> 45505
    @{ Row=22; Col=3; Value=331 }  # C22: 328 -> 331
    @{ Row=22; Col=4; Value=375 }  # D22: 373 -> 375
    @{ Row=22; Col=6; Value=45444 }  # F22: 45352 -> 45444
    @{ Row=22; Col=8; Value=45505 }  # H22: 45446 -> 45505
    @{ Row=23; Col=4; Value=323 }  # D23: 321 -> 323
    @{ Row=23; Col=5; Value=30225 }  # E23: 30133 -> 30225
    @{ Row=23; Col=6; Value=45413 }  # F23: 45323 -> 45413
    @{ Row=23; Col=8; Value=45505 }  # H23: 45446 -> 45505
    @{ Row=24; Col=3; Value=341 }  # C24: 340 -> 341
    @{ Row=24; Col=4; Value=321 }  # D24: 319 -> 321
    @{ Row=24; Col=6; Value=45444 }  # F24: 45413 -> 45444
    @{ Row=24; Col=8; Value=45505 }  # H24: 45446 -> 45505
    @{ Row=25; Col=5; Value=30225 }  # E25: 30164 -> 30225
    @{ Row=25; Col=6; Value=45413 }  # F25: 45352 -> 45413
    @{ Row=25; Col=7; Value=30319 }  # G25: 30256 -> 30319
    @{ Row=25; Col=8; Value=45505 }  # H25: 45446 -> 45505
    @{ Row=26; Col=4; Value=388 }  # D26: 386 -> 388
    @{ Row=26; Col=5; Value=30256 }  # E26: 30195 -> 30256
    @{ Row=26; Col=6; Value=45444 }  # F26: 45383 -> 45444
    @{ Row=26; Col=8; Value=45505 }  # H26: 45446 -> 45505
    @{ Row=27; Col=4; Value=225 }  # D27: 223 -> 225
    @{ Row=27; Col=5; Value=30225 }  # E27: 30195 -> 30225
    @{ Row=27; Col=6; Value=45413 }  # F27: 45383 -> 45413
    @{ Row=27; Col=8; Value=45505 }  # H27: 45446 -> 45505
    @{ Row=28; Col=3; Value=402 }  # C28: 399 -> 402
    @{ Row=28; Col=6; Value=45444 }  # F28: 45352 -> 45444
    @{ Row=28; Col=7; Value=30319 }  # G28: 30256 -> 30319
    @{ Row=28; Col=8; Value=45505 }  # H28: 45446 -> 45505
    @{ Row=29; Col=7; Value=30319 }  # G29: 30256 -> 30319
    @{ Row=29; Col=8; Value=45505 }  # H29: 45446 -> 45505
    @{ Row=30; Col=3; Value=462 }  # C30: 460 -> 462
    @{ Row=30; Col=4; Value=406 }  # D30: 404 -> 406
    @{ Row=30; Col=6; Value=45444 }  # F30: 45383 -> 45444
    @{ Row=30; Col=8; Value=45505 }  # H30: 45446 -> 45505
    @{ Row=31; Col=3; Value=412 }  # C31: 411 -> 412
    @{ Row=31; Col=4; Value=325 }  # D31: 323 -> 325
    @{ Row=31; Col=6; Value=45413 }  # F31: 45383 -> 45413
    @{ Row=31; Col=8; Value=45505 }  # H31: 45446 -> 45505
    @{ Row=32; Col=4; Value=418 }  # D32: 416 -> 418
    @{ Row=32; Col=5; Value=30256 }  # E32: 30225 -> 30256
    @{ Row=32; Col=6; Value=45444 }  # F32: 45413 -> 45444
    @{ Row=32; Col=8; Value=45505 }  # H32: 45446 -> 45505
    @{ Row=33; Col=3; Value=474 }  # C33: 472 -> 474
    @{ Row=33; Col=4; Value=325 }  # D33: 323 -> 325
    @{ Row=33; Col=6; Value=45444 }  # F33: 45383 -> 45444
    @{ Row=33; Col=8; Value=45505 }  # H33: 45446 -> 45505
    @{ Row=34; Col=3; Value=367 }  # C34: 365 -> 367
    @{ Row=34; Col=4; Value=375 }  # D34: 373 -> 375
    @{ Row=34; Col=6; Value=45444 }  # F34: 45383 -> 45444
    @{ Row=34; Col=8; Value=45505 }  # H34: 45446 -> 45505
)

foreach ($update in $cellUpdates) {
    $ws.Cells.Item($update.Row, $update.Col).Value = $update.Value
}
